# Adding tech, fuel, and storage to filter
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Fuel_selection: add a new fuel "Heat_Low_DistrictHeat" as a new row
#    inserted right after "Cool_Low_Building" (before "Biofuel").
# ---------------------------------------------------------------------
$wsFuel = $wb.Worksheets.Item("Fuel_selection")
$wsFuel.Rows.Item(21).Insert()
$wsFuel.Range("A21").Value = "Heat_Low_DistrictHeat"
$wsFuel.Range("B21").Value = 1

# widen column A to fit the longer fuel names
$wsFuel.Columns.Item(1).ColumnWidth = 24.1

$wsFuel.Range("E15").Select()

# ---------------------------------------------------------------------
# 2) Storage_selection: add two new storage types at the bottom of the
#    list: "S_Heat_HLB" and "S_Heat_HLDH".
# ---------------------------------------------------------------------
$wsStorage = $wb.Worksheets.Item("Storage_selection")
$wsStorage.Range("A10").Value = "S_Heat_HLB"
$wsStorage.Range("B10").Value = 1
$wsStorage.Range("A11").Value = "S_Heat_HLDH"
$wsStorage.Range("B11").Value = 1

# set explicit column widths for the storage sheet
$wsStorage.Columns.Item(1).ColumnWidth = 20.1
$wsStorage.Columns.Item(2).ColumnWidth = 17.3

$wsStorage.Range("A12").Select()

# ---------------------------------------------------------------------
# 3) Technology_selection: add four new heat-pump technologies at the
#    bottom of the list.
# ---------------------------------------------------------------------
$wsTech = $wb.Worksheets.Item("Technology_selection")
$wsTech.Range("A163").Value = "HLDH_Heatpump_Air"
$wsTech.Range("B163").Value = 1
$wsTech.Range("A164").Value = "HLDH_Heatpump_ExcessHeat"
$wsTech.Range("B164").Value = 1
$wsTech.Range("A165").Value = "HLI_Heatpump"
$wsTech.Range("B165").Value = 1
$wsTech.Range("A166").Value = "HMLI_Heatpump"
$wsTech.Range("B166").Value = 1

# Technology_selection ends up as the active / selected sheet
$wsTech.Activate()
$wsTech.Range("A161").Select()
